# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Typhon_Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 28000
$ws.Range("J63").Value = 28000
$ws.Range("L63").Value = 28000
$ws.Range("N63").Value = -29248

$ws.Range("H66").Value = 28000
$ws.Range("J66").Value = 28000
$ws.Range("L66").Value = 84000
$ws.Range("N66").Value = -90240

$ws.Range("H96").Value = 41667064
$ws.Range("I96").Value = 41667064
$ws.Range("K96").Value = 125001192
$ws.Range("M96").Value = -124999819

$ws.Range("H132").Value = 4462.9473
$ws.Range("I132").Value = 4655.3335
$ws.Range("K132").Value = 13966.0005
$ws.Range("M132").Value = -11436.0005

$ws.Range("H135").Value = 26324558
$ws.Range("I135").Value = 829.3571
$ws.Range("K135").Value = 7464.2139
$ws.Range("M135").Value = -4929.2139

$ws.Range("H138").Value = 2317.797
$ws.Range("I138").Value = 2590.6667
$ws.Range("J138").Value = 2260.3508
$ws.Range("K138").Value = 7772.000100000001
$ws.Range("L138").Value = 6781.0524
$ws.Range("M138").Value = -2632.000100000001
$ws.Range("N138").Value = -17061.0524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1666.8387
$ws.Range("I2").Value = 1577.5
$ws.Range("J2").Value = 1973.1428
$ws.Range("K2").Value = 1577.5
$ws.Range("L2").Value = 1973.1428
$ws.Range("M2").Value = -1464.5
$ws.Range("N2").Value = -2199.1428

$ws.Range("H32").Value = 5106.3125
$ws.Range("I32").Value = 4469.9785
$ws.Range("K32").Value = 4469.9785
$ws.Range("M32").Value = -4182.9785

$ws.Range("H74").Value = 76924104
$ws.Range("I74").Value = 90909760
$ws.Range("K74").Value = 90909760
$ws.Range("M74").Value = -90908886

$ws.Range("H77").Value = 76924104
$ws.Range("I77").Value = 90909760
$ws.Range("K77").Value = 454548800
$ws.Range("M77").Value = -454544432

$ws.Range("H88").Value = 201910.8
$ws.Range("J88").Value = 335518
$ws.Range("L88").Value = 335518
$ws.Range("N88").Value = -336330

$ws.Range("H91").Value = 201910.8
$ws.Range("J91").Value = 335518
$ws.Range("L91").Value = 335518
$ws.Range("N91").Value = -338326

$ws.Range("H116").Value = 1666.8387
$ws.Range("I116").Value = 1577.5
$ws.Range("J116").Value = 1973.1428
$ws.Range("K116").Value = 1577.5
$ws.Range("L116").Value = 1973.1428
$ws.Range("M116").Value = 716.5
$ws.Range("N116").Value = -6561.1428

$ws.Range("H132").Value = 18745.9
$ws.Range("I132").Value = 2093.95
$ws.Range("J132").Value = 52049.8
$ws.Range("K132").Value = 6281.849999999999
$ws.Range("L132").Value = 156149.4
$ws.Range("M132").Value = -3751.849999999999
$ws.Range("N132").Value = -161209.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1666.8387
$ws.Range("I3").Value = 1577.5
$ws.Range("J3").Value = 1973.1428
$ws.Range("K3").Value = 1577.5
$ws.Range("L3").Value = 1973.1428
$ws.Range("M3").Value = -1463.5
$ws.Range("N3").Value = -2201.1428

$ws.Range("H86").Value = 1841.75
$ws.Range("I86").Value = 1457.0454
$ws.Range("K86").Value = 1457.0454
$ws.Range("M86").Value = -334.0454

$ws.Range("H89").Value = 1841.75
$ws.Range("I89").Value = 1457.0454
$ws.Range("K89").Value = 7285.227
$ws.Range("M89").Value = -1669.227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15722.074
$ws.Range("I31").Value = 26523.23
$ws.Range("K31").Value = 26523.23
$ws.Range("M31").Value = -26228.23

$ws.Range("H34").Value = 15722.074
$ws.Range("I34").Value = 26523.23
$ws.Range("K34").Value = 26523.23
$ws.Range("M34").Value = -26321.23

$ws.Range("H56").Value = 11240
$ws.Range("I56").Value = 2480
$ws.Range("J56").Value = 20000
$ws.Range("K56").Value = 2480
$ws.Range("L56").Value = 20000
$ws.Range("M56").Value = -1635
$ws.Range("N56").Value = -21690

$ws.Range("H58").Value = 34746.133
$ws.Range("I58").Value = 1580.091
$ws.Range("J58").Value = 125952.75
$ws.Range("K58").Value = 1580.091
$ws.Range("L58").Value = 125952.75
$ws.Range("M58").Value = -1377.091
$ws.Range("N58").Value = -126358.75

$ws.Range("H86").Value = 12125.066
$ws.Range("I86").Value = 3255.5454
$ws.Range("J86").Value = 36516.25
$ws.Range("K86").Value = 3255.5454
$ws.Range("L86").Value = 36516.25
$ws.Range("M86").Value = -2132.5454
$ws.Range("N86").Value = -38762.25

$ws.Range("H89").Value = 12125.066
$ws.Range("I89").Value = 3255.5454
$ws.Range("J89").Value = 36516.25
$ws.Range("K89").Value = 16277.727
$ws.Range("L89").Value = 182581.25
$ws.Range("M89").Value = -10661.727
$ws.Range("N89").Value = -193813.25

$ws.Range("H92").Value = 49899
$ws.Range("J92").Value = 49899
$ws.Range("L92").Value = 49899
$ws.Range("N92").Value = -54891

$ws.Range("H136").Value = 34746.133
$ws.Range("I136").Value = 1580.091
$ws.Range("J136").Value = 125952.75
$ws.Range("K136").Value = 4740.272999999999
$ws.Range("L136").Value = 377858.25
$ws.Range("M136").Value = -2190.272999999999
$ws.Range("N136").Value = -382958.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 756.53845
$ws.Range("I113").Value = 737.6667
$ws.Range("J113").Value = 762.2
$ws.Range("K113").Value = 2213.0001
$ws.Range("L113").Value = 2286.6
$ws.Range("M113").Value = -43.0001000000002
$ws.Range("N113").Value = -6626.6

$ws.Range("H131").Value = 718.07
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 721.28284
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 2163.84852
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -12243.84852

$ws.Range("H134").Value = 3938
$ws.Range("I134").Value = 1435
$ws.Range("J134").Value = 6441
$ws.Range("K134").Value = 4305
$ws.Range("L134").Value = 19323
$ws.Range("M134").Value = 765
$ws.Range("N134").Value = -29463

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2411552
$ws.Range("I70").Value = 9355.723
$ws.Range("J70").Value = 7816493.5
$ws.Range("K70").Value = 9355.723
$ws.Range("L70").Value = 7816493.5
$ws.Range("M70").Value = -9085.723
$ws.Range("N70").Value = -7817033.5

$ws.Range("H73").Value = 2411552
$ws.Range("I73").Value = 9355.723
$ws.Range("J73").Value = 7816493.5
$ws.Range("K73").Value = 9355.723
$ws.Range("L73").Value = 7816493.5
$ws.Range("M73").Value = -8419.723
$ws.Range("N73").Value = -7818365.5

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H80").Value = 3565.75
$ws.Range("I80").Value = 3257.2727
$ws.Range("K80").Value = 3257.2727
$ws.Range("M80").Value = -2259.2727

$ws.Range("H83").Value = 3565.75
$ws.Range("I83").Value = 3257.2727
$ws.Range("K83").Value = 16286.3635
$ws.Range("M83").Value = -11294.3635

$ws.Range("H113").Value = 3115.3845
$ws.Range("J113").Value = 3780
$ws.Range("L113").Value = 3780
$ws.Range("N113").Value = -8120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 34000
$ws.Range("J64").Value = 34000
$ws.Range("L64").Value = 34000
$ws.Range("N64").Value = -34450

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 34000
$ws.Range("J67").Value = 34000
$ws.Range("L67").Value = 34000
$ws.Range("N67").Value = -35560

$ws.Range("H93").Value = 2869.4
$ws.Range("I93").Value = 2521.111
$ws.Range("J93").Value = 6004
$ws.Range("K93").Value = 2521.111
$ws.Range("L93").Value = 6004
$ws.Range("M93").Value = -1273.111
$ws.Range("N93").Value = -8500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1593.6
$ws.Range("I122").Value = 1629.9048
$ws.Range("J122").Value = 1403
$ws.Range("K122").Value = 4889.7144
$ws.Range("L122").Value = 4209
$ws.Range("M122").Value = -2439.7144
$ws.Range("N122").Value = -9109

$ws.Range("H136").Value = 33335086
$ws.Range("I136").Value = 40001388
$ws.Range("J136").Value = 3580.8
$ws.Range("K136").Value = 120004164
$ws.Range("L136").Value = 10742.4
$ws.Range("M136").Value = -120001614
$ws.Range("N136").Value = -15842.4

$ws.Range("H139").Value = 52715
$ws.Range("J139").Value = 52715
$ws.Range("L139").Value = 52715
$ws.Range("N139").Value = -62995
